# Applies the BDbDT workbook update:
#  - BDbDT!B2:AF9 literal zero values are replaced with SUMIFS() formulas
#    that total TOTAL_DEATHS from the 'Census T3' sheet by RACE_HISP (col A),
#    SEX (col B) and YEAR (col C$1 of BDbDT), with two rows needing an extra
#    adjustment term.
#  - The active / selected worksheet changes from "BDbDT" back to "About".

$wb = $excel.ActiveWorkbook
$wsBDbDT = $wb.Worksheets.Item("BDbDT")
$wsCensus = $wb.Worksheets.Item("Census T3")

# Column letters for the year columns B (2020) through AF (2050).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF")

# Row -> (RACE_HISP criteria, SEX criteria, extra formula suffix template)
# The suffix template, if present, is appended to the base SUMIFS() call;
# {0} is replaced with the current column letter.
$rows = @(
    @{ Row = 2; Race = 0; Sex = 1; Suffix = "" },
    @{ Row = 3; Race = 0; Sex = 2; Suffix = "" },
    @{ Row = 4; Race = 1; Sex = 0; Suffix = "" },
    @{ Row = 5; Race = 2; Sex = 0; Suffix = "" },
    @{ Row = 6; Race = 4; Sex = 0; Suffix = "" },
    @{ Row = 7; Race = 0; Sex = 0; Suffix = "-SUM({0}4:{0}6)" },
    @{ Row = 8; Race = 8; Sex = 0; Suffix = "" },
    @{ Row = 9; Race = 0; Sex = 0; Suffix = "-{0}8" }
)

foreach ($rowInfo in $rows) {
    $r = $rowInfo.Row
    $race = $rowInfo.Race
    $sex = $rowInfo.Sex
    $suffixTemplate = $rowInfo.Suffix

    foreach ($c in $cols) {
        $base = "SUMIFS('Census T3'!`$D:`$D,'Census T3'!`$A:`$A,$race,'Census T3'!`$B:`$B,$sex,'Census T3'!`$C:`$C,$c`$1)"
        if ($suffixTemplate -ne "") {
            $suffix = $suffixTemplate -f $c
            $formula = "=" + $base + $suffix
        } else {
            $formula = "=" + $base
        }
        $wsBDbDT.Range("$c$r").Formula = $formula
    }
}

# Restore "About" as the active / selected sheet (it was the active tab
# before BDbDT was last worked on), which also clears BDbDT's stale
# tabSelected/topLeftCell/selection view state.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()

Write-Host "BDbDT formulas applied; active sheet set to About."
